$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2: C2 "asd" -> "Merchant" (gets the "alt/plain" look used for the new
# merchant/product grouping cells throughout the sheet)
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Merchant"
$ws.Range("C2").Font.Name = "Calibri"
$ws.Range("C2").Font.Size = 11

# ---------------------------------------------------------------------------
# Row 4: B4 "Finca" -> "Otros animales domesticos" (same alt look)
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "Otros animales domesticos"
$ws.Range("B4").Font.Name = "Calibri"
$ws.Range("B4").Font.Size = 11

# ---------------------------------------------------------------------------
# Row 5: B5 "Documentacion" -> "Alimentos "; C5 stops being a shared formula
# (=C4) and becomes the literal value "Product" with the alt look; D5:E5
# becomes its own shared formula block referencing D4 instead of C4.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "Alimentos "
$ws.Range("C5").Formula = ""
$ws.Range("C5").Value = "Product"
$ws.Range("C5").Font.Name = "Calibri"
$ws.Range("C5").Font.Size = 11
$ws.Range("D5").Formula = "=D4"
$ws.Range("E5").Formula = "=D4"

# ---------------------------------------------------------------------------
# Rows 6-9: rename categories (values only, formulas untouched)
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = "Aseo y belleza animal"
$ws.Range("B7").Value = "Servicios veterinarios"
$ws.Range("B8").Value = "Juguetes y accesorios"
$ws.Range("B9").Value = "Servicios especializados"
$ws.Range("B9").Font.Name = "Calibri"
$ws.Range("B9").Font.Size = 11

# ---------------------------------------------------------------------------
# New rows 10 and 11
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Destacados"
$ws.Range("C10").Formula = "=C9"
$ws.Range("D10").Formula = "=C9"
$ws.Range("E10").Formula = "=C9"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "FAQs"
$ws.Range("C11").Value = "Article"
$ws.Range("D11").Formula = "=D10"
$ws.Range("E11").Formula = "=D10"

# ---------------------------------------------------------------------------
# New rows 12-16: label in column B (alt look, last one plain) plus two
# placeholder formatted-but-empty cells in F/G (italic grey note column and
# a grey right-aligned percent/number column).
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = "Documentacion"
$ws.Range("B12").Font.Name = "Calibri"
$ws.Range("B12").Font.Size = 11

$ws.Range("B13").Value = "Vacunas"
$ws.Range("B13").Font.Name = "Calibri"
$ws.Range("B13").Font.Size = 11

$ws.Range("B14").Value = "Medicinas"
$ws.Range("B14").Font.Name = "Calibri"
$ws.Range("B14").Font.Size = 11

$ws.Range("B15").Value = "Comida"
$ws.Range("B15").Font.Name = "Calibri"
$ws.Range("B15").Font.Size = 11

$ws.Range("B16").Value = "Juguetes"

$ws.Range("B17").Value = "Destacados"

$ws.Range("F12:F16").Font.Italic = $true
$ws.Range("F12:F16").Font.Size = 11
$ws.Range("F12:F16").Font.Color = 10066329
$ws.Range("F12:F16").Font.Name = "Calibri"

$ws.Range("G12:G16").Font.Size = 11
$ws.Range("G12:G16").Font.Color = 10066329
$ws.Range("G12:G16").Font.Name = "Calibri"
$ws.Range("G12:G16").HorizontalAlignment = -4152

$ws.Range("G12:G14").NumberFormat = "0%"

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 24.2
$ws.Columns.Item(6).ColumnWidth = 19.6
